# Quiz_Graded.xlsx edit: new grading API adapter.
# The "Reason N" explanation columns are dropped; the "Soru N Puan" score
# columns are renumbered/shifted left to fill the gap (AD122 -> Y122).
# Row 2 is re-graded under the new schema; rows 3-17 (previously graded
# under the old schema) are cleared so they can be re-graded; rows 18+
# were already ungraded and stay that way.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Update the header labels for the columns that survive the shift.
$ws.Range("V1").Value = "Soru 2 Puan"
$ws.Range("W1").Value = "Soru 3 Puan"
$ws.Range("X1").Value = "Soru 4 Puan"
$ws.Range("Y1").Value = "Soru 5 Puan"

# 2) Re-graded row (student 20240808081): new scores under the new schema.
$ws.Range("T2").Value = 95
$ws.Range("U2").Value = 15
$ws.Range("V2").Value = 20
$ws.Range("W2").Value = 20
$ws.Range("X2").Value = 20
$ws.Range("Y2").Value = 20

# 3) Clear the stale grading data (old schema) for rows 3-17 so the
#    new grading pass can fill them in later.
$ws.Range("T3:Y17").ClearContents()

# 4) Drop the now-obsolete "Reason" columns entirely (Z:AD), shifting
#    nothing else left of them and shrinking the sheet dimension to
#    A1:Y122.
$ws.Range("Z1:AD122").EntireColumn.Delete()
